$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 44, pushing existing rows 44-145 down to 46-147.
$ws.Rows("44:45").Insert()

# Populate the two newly inserted rows with the new data entries.
$ws.Range("A44").Value = 7
$ws.Range("B44").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C44").Value = "Ñuble"
$ws.Range("D44").Value = 44469
$ws.Range("E44").Value = 16
$ws.Range("F44").Value = 100112002
$ws.Range("G44").Value = "Pimiento"
$ws.Range("H44").Value = "Zafiro rojo"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 120
$ws.Range("K44").Value = 42000
$ws.Range("L44").Value = 43000
$ws.Range("M44").Value = 42500
$ws.Range("N44").Value = "$/caja 15 kilos"
$ws.Range("O44").Value = "Región de Arica y Parinacota"
$ws.Range("P44").Value = 2833
$ws.Range("Q44").Value = 15
$ws.Range("R44").Value = "Hortaliza"

$ws.Range("A45").Value = 7
$ws.Range("B45").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C45").Value = "Ñuble"
$ws.Range("D45").Value = 44469
$ws.Range("E45").Value = 16
$ws.Range("F45").Value = 100112002
$ws.Range("G45").Value = "Pimiento"
$ws.Range("H45").Value = "Zafiro verde"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 120
$ws.Range("K45").Value = 35000
$ws.Range("L45").Value = 36000
$ws.Range("M45").Value = 35500
$ws.Range("N45").Value = "$/caja 15 kilos"
$ws.Range("O45").Value = "Región de Arica y Parinacota"
$ws.Range("P45").Value = 2367
$ws.Range("Q45").Value = 15
$ws.Range("R45").Value = "Hortaliza"
